$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-10 03:18:19'
$ws.Range('G2').Value = '207 cm'
$ws.Range('I2').Value = '7.0 mm'
$ws.Range('M2').Value = '0.3 °C 2:29 TU'
$ws.Range('N2').Value = '-0.4 °C 1:30 TU'
$ws.Range('O2').Value = '-0.1 °C'
$ws.Range('E3').Value = '2026-02-10 03:18:22'
$ws.Range('I3').Value = '5.0 mm'
$ws.Range('E4').Value = '2026-02-10 03:18:24'
$ws.Range('H4').NumberFormat = '@'
$ws.Range('H4').Value = '86%'
$ws.Range('O4').Value = '9.1 °C'
$ws.Range('E5').Value = '2026-02-10 03:18:26'
$ws.Range('G5').Value = '129 cm'
$ws.Range('I5').Value = '7.3 mm'
$ws.Range('N5').Value = '-1.1 °C 2:33 TU'
$ws.Range('O5').Value = '-0.5 °C'
$ws.Range('E6').Value = '2026-02-10 03:18:29'
$ws.Range('L6').Value = '12.6 km/h - 335º 2:48 TU'
$ws.Range('E7').Value = '2026-02-10 03:18:31'
$ws.Range('J7').Value = '1005.3 hPa'
$ws.Range('E8').Value = '2026-02-10 03:18:33'
$ws.Range('J8').Value = '1005.4 hPa'
$ws.Range('E9').Value = '2026-02-10 03:18:36'
$ws.Range('L9').Value = '6.5 km/h - 273º 2:54 TU'
$ws.Range('N9').Value = '6.5 °C 2:47 TU'
$ws.Range('O9').Value = '7.2 °C'
$ws.Range('E10').Value = '2026-02-10 03:18:38'
$ws.Range('N10').Value = '7.1 °C 2:48 TU'
$ws.Range('E11').Value = '2026-02-10 03:18:40'
$ws.Range('I11').Value = '0.1 mm'
$ws.Range('N11').Value = '2.6 °C 2:58 TU'
$ws.Range('O11').Value = '3.2 °C'
$ws.Range('E12').Value = '2026-02-10 03:18:43'
$ws.Range('N12').Value = '7.0 °C 2:44 TU'
$ws.Range('E13').Value = '2026-02-10 03:18:45'
$ws.Range('I13').Value = '0.8 mm'
$ws.Range('N13').Value = '2.5 °C 2:56 TU'
$ws.Range('E14').Value = '2026-02-10 03:18:47'
$ws.Range('H14').NumberFormat = '@'
$ws.Range('H14').Value = '99%'
$ws.Range('N14').Value = '9.1 °C 2:44 TU'
$ws.Range('O14').Value = '10.0 °C'
$ws.Range('E15').Value = '2026-02-10 03:18:50'
$ws.Range('N15').Value = '6.6 °C 2:59 TU'
$ws.Range('O15').Value = '7.6 °C'
$ws.Range('E16').Value = '2026-02-10 03:18:52'
$ws.Range('I16').Value = '5.1 mm'
$ws.Range('O16').Value = '-0.5 °C'
$ws.Range('E17').Value = '2026-02-10 03:18:54'
$ws.Range('G17').Value = '1 cm'
$ws.Range('M17').Value = '2.2 °C 2:59 TU'
$ws.Range('O17').Value = '1.4 °C'
$ws.Range('E18').Value = '2026-02-10 03:18:57'
$ws.Range('H18').NumberFormat = '@'
$ws.Range('H18').Value = '96%'
$ws.Range('N18').Value = '7.3 °C 2:31 TU'
$ws.Range('O18').Value = '8.1 °C'
$ws.Range('E19').Value = '2026-02-10 03:18:59'
$ws.Range('E20').Value = '2026-02-10 03:19:02'
$ws.Range('I20').Value = '1.5 mm'
$ws.Range('M20').Value = '-1.0 °C 2:59 TU'
$ws.Range('O20').Value = '-1.6 °C'
$ws.Range('E21').Value = '2026-02-10 03:19:04'
$ws.Range('I21').Value = '2.1 mm'
$ws.Range('N21').Value = '4.2 °C 2:57 TU'
$ws.Range('O21').Value = '4.4 °C'
$ws.Range('E22').Value = '2026-02-10 03:19:07'
$ws.Range('M22').Value = '-1.6 °C 2:47 TU'
$ws.Range('O22').Value = '-2.2 °C'
$ws.Range('E23').Value = '2026-02-10 03:19:09'
$ws.Range('G23').Value = '179 cm'
$ws.Range('H23').NumberFormat = '@'
$ws.Range('H23').Value = '91%'
$ws.Range('I23').Value = '5.6 mm'
$ws.Range('O23').Value = '-0.7 °C'
$ws.Range('E24').Value = '2026-02-10 03:19:11'
$ws.Range('L24').Value = '10.1 km/h - 328º 2:58 TU'
$ws.Range('E25').Value = '2026-02-10 03:19:14'
$ws.Range('G25').Value = '114 cm'
$ws.Range('I25').Value = '2.1 mm'
$ws.Range('L25').Value = '28.4 km/h - 301º 2:47 TU'
$ws.Range('M25').Value = '1.1 °C 2:55 TU'
$ws.Range('O25').Value = '-1.1 °C'
$ws.Range('E26').Value = '2026-02-10 03:19:16'
$ws.Range('J26').Value = '1005.1 hPa'
$ws.Range('N26').Value = '2.3 °C 2:57 TU'
$ws.Range('O26').Value = '3.0 °C'
$ws.Range('E27').Value = '2026-02-10 03:19:19'
$ws.Range('I27').Value = '1.2 mm'
$ws.Range('L27').Value = '32.4 km/h - 255º 2:53 TU'
$ws.Range('O27').Value = '-0.8 °C'
$ws.Range('E28').Value = '2026-02-10 03:19:21'
$ws.Range('J28').Value = '1005.5 hPa'
$ws.Range('N28').Value = '4.9 °C 2:55 TU'
$ws.Range('O28').Value = '5.7 °C'
$ws.Range('E29').Value = '2026-02-10 03:19:24'
$ws.Range('M29').Value = '10.6 °C 2:34 TU'
$ws.Range('O29').Value = '9.6 °C'
$ws.Range('E30').Value = '2026-02-10 03:19:26'
$ws.Range('N30').Value = '7.1 °C 2:45 TU'
$ws.Range('O30').Value = '7.6 °C'
$ws.Range('E31').Value = '2026-02-10 03:19:28'
$ws.Range('N31').Value = '8.9 °C 2:59 TU'
$ws.Range('O31').Value = '9.2 °C'
$ws.Range('E32').Value = '2026-02-10 03:19:31'
$ws.Range('I32').Value = '0.6 mm'
$ws.Range('L32').Value = '22.3 km/h - 316º 2:57 TU'
$ws.Range('M32').Value = '8.1 °C 2:49 TU'
$ws.Range('O32').Value = '7.3 °C'
$ws.Range('E33').Value = '2026-02-10 03:19:33'
$ws.Range('I33').Value = '0.7 mm'
$ws.Range('N33').Value = '1.8 °C 2:37 TU'
$ws.Range('E34').Value = '2026-02-10 03:19:36'
$ws.Range('I34').Value = '0.9 mm'
$ws.Range('O34').Value = '2.3 °C'
$ws.Range('E35').Value = '2026-02-10 03:19:38'
$ws.Range('L35').Value = '62.6 km/h - 271º 2:38 TU'
$ws.Range('O35').Value = '10.6 °C'
$ws.Range('E36').Value = '2026-02-10 03:19:41'
$ws.Range('H36').NumberFormat = '@'
$ws.Range('H36').Value = '97%'
$ws.Range('L36').Value = '31.7 km/h - 11º 2:49 TU'
$ws.Range('O36').Value = '8.8 °C'
$ws.Range('E37').Value = '2026-02-10 03:19:43'
$ws.Range('H37').NumberFormat = '@'
$ws.Range('H37').Value = '94%'
$ws.Range('N37').Value = '3.8 °C 2:59 TU'
$ws.Range('O37').Value = '4.5 °C'
$ws.Range('E38').Value = '2026-02-10 03:19:45'
$ws.Range('E39').Value = '2026-02-10 03:19:48'
$ws.Range('I39').Value = '0.7 mm'
$ws.Range('L39').Value = '43.9 km/h - 324º 2:53 TU'
$ws.Range('E40').Value = '2026-02-10 03:19:50'
$ws.Range('I40').Value = '1.4 mm'
$ws.Range('N40').Value = '4.4 °C 2:55 TU'
$ws.Range('O40').Value = '5.0 °C'
$ws.Range('E41').Value = '2026-02-10 03:19:52'
$ws.Range('H41').NumberFormat = '@'
$ws.Range('H41').Value = '97%'
$ws.Range('N41').Value = '9.0 °C 2:51 TU'
$ws.Range('O41').Value = '10.6 °C'
$ws.Range('E42').Value = '2026-02-10 03:19:55'
$ws.Range('E43').Value = '2026-02-10 03:19:57'
$ws.Range('E44').Value = '2026-02-10 03:19:59'
$ws.Range('I44').Value = '5.0 mm'
$ws.Range('O44').Value = '-0.6 °C'
$ws.Range('E45').Value = '2026-02-10 03:20:02'
$ws.Range('H45').NumberFormat = '@'
$ws.Range('H45').Value = '98%'
$ws.Range('I45').Value = '7.5 mm'
$ws.Range('E46').Value = '2026-02-10 03:20:04'
$ws.Range('N46').Value = '9.1 °C 2:31 TU'
